$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D4, E4, G4 and F4 used a "Noto Sans CJK SC" font variant without a charset
# declaration; the edit re-points them at the (already present) variant that
# does declare one. F4 additionally keeps its wrap-text formatting, so set
# that before the font so the engine folds it into the same style record.
$ws.Range("D4").Font.Name = "Noto Sans CJK SC"
$ws.Range("E4").Font.Name = "Noto Sans CJK SC"
$ws.Range("G4").Font.Name = "Noto Sans CJK SC"
$ws.Range("F4").WrapText = $true
$ws.Range("F4").Font.Name = "Noto Sans CJK SC"

# New remark row appended below the existing data.
$c = $ws.Range("F5")
$c.Value = "dofも初期位置ランダム要素を加えてみたらどうなるだろうか？"

# "dof" keeps the Latin/Arial run, the rest of the sentence uses the CJK font
# -- matching the existing mixed-run comments elsewhere in the sheet.
$c.Characters(1, 3).Font.Name = "Arial"
$c.Characters(4, 27).Font.Name = "Noto Sans CJK SC"

# Selection moved on to the next empty cell after the edit.
$ws.Range("F6").Select() | Out-Null
